$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying file list was regenerated: two pairs of rows (sharing the
# same "Review date") swapped order, and one file was renamed as part of
# the regeneration.

# Group with Review date 2025-10-20: move "Procedures/Inadvertent Catheter
# Placement Guideline.pdf" above "Policies_and_admin/Pet Visitation.pdf"
$ws.Range("A163").Value = "Procedures/Inadvertent Catheter Placement Guideline.pdf"
$ws.Range("A164").Value = "Policies_and_admin/Pet Visitation.pdf"

# Group with Review date 2026-01-20: rename "Neurological/Status guideline
# 24.pdf" to "Neurological/Status epilepticus.pdf" and move it above
# "Drugs/Thiopentone.pdf"
$ws.Range("A170").Value = "Neurological/Status epilepticus.pdf"
$ws.Range("A171").Value = "Drugs/Thiopentone.pdf"
